# Generate Report for Handoff
#
# The previous handoff round (315291f9-... / 64e4bdde-...) has been
# superseded by a new handoff round (b0489487-... / ffff7664ed3f-...):
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - New source/target file names + timestamps
#   - The "Latest Target File" / "Latest Handback File" columns no longer
#     apply (nothing has been handed back yet), so those two columns of
#     data are cleared out on the per-language sheets.

$wb = $excel.ActiveWorkbook

$oldMd1 = "315291f9-87ee-4e64-bf67-ad229cc2873a.md"
$oldMd2 = "64e4bdde-47cc-465d-afe3-3f7921c6a394.md"
$newMd1 = "b0489487-19e5-4b3f-87d8-aaa577d701b1.md"
$newMd2 = "ffff7664ed3f-a6c2-434a-9dd7-0b212984201f.md"

$newStatus = "Ready for handoff"
$newHandoffDate = "2016-03-25 03:24:57"

$newZhXlf = "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf"
$newDeXlf = "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf"
$newZhHandoffDatetime = "2016-03-25 03:24:53"
$newDeHandoffDatetime = "2016-03-25 03:24:57"
$clearedHandbackDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMd1
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffDatetime
$wsZh.Range("H2").Value = $clearedHandbackDatetime

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $newZhXlf
$wsZh.Range("E3").Value = $newZhHandoffDatetime
$wsZh.Range("H3").Value = $clearedHandbackDatetime

# Drop the now-irrelevant "Latest Target File" / "Latest Handback File"
# columns (F, G) for both data rows.
$zhRemove = @("`$F`$2", "`$G`$2", "`$F`$3", "`$G`$3")
foreach ($target in $zhRemove) {
    foreach ($h in $wsZh.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            $h.Delete()
            break
        }
    }
}
$wsZh.Range("F2:G2").Clear()
$wsZh.Range("F3:G3").Clear()

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMd1
    } elseif ($addr -eq "`$D`$2") {
        $h.TextToDisplay = $newZhXlf
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = $newMd2
    } elseif ($addr -eq "`$D`$3") {
        $h.TextToDisplay = $newZhXlf
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeHandoffDatetime
$wsDe.Range("H2").Value = $clearedHandbackDatetime

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $newDeXlf
$wsDe.Range("E3").Value = $newDeHandoffDatetime
$wsDe.Range("H3").Value = $clearedHandbackDatetime

$deRemove = @("`$F`$2", "`$G`$2", "`$F`$3", "`$G`$3")
foreach ($target in $deRemove) {
    foreach ($h in $wsDe.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            $h.Delete()
            break
        }
    }
}
$wsDe.Range("F2:G2").Clear()
$wsDe.Range("F3:G3").Clear()

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = $newMd1
    } elseif ($addr -eq "`$D`$2") {
        $h.TextToDisplay = $newDeXlf
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = $newMd2
    } elseif ($addr -eq "`$D`$3") {
        $h.TextToDisplay = $newDeXlf
    }
}

Write-Host "Handoff report regenerated."
